$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 73-79 (answer judged "tie") get a new E column value of 0
$tieRows = 73..79
foreach ($r in $tieRows) {
    $ws.Cells.Item($r, 5).Value = 0
}

# Rows 80-81 (answer judged in favour of the other side) get -1
$ws.Cells.Item(80, 5).Value = -1
$ws.Cells.Item(81, 5).Value = -1

# Reflect where the author scrolled to / selected before saving
$ws.Range("E82").Select()
$excel.ActiveWindow.ScrollRow = 81
$excel.ActiveWindow.ScrollColumn = 1
